# HighLevelSequenceDiagrams.pptx update (UG/DG diagrams)
# Converts EMU -> points (1 pt = 12700 EMU) since the Shape position/size
# COM properties (Left/Top/Width/Height) are expressed in points.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function EMU([double]$v) { return $v / 12700.0 }

# ---------------------------------------------------------------------
# 1) Refresh the "02-Apr-19" datetimeFigureOut fields (slide master, all
#    slide layouts and the notes master) to "4/15/2019".
# ---------------------------------------------------------------------
$sm = $p.SlideMaster
$sm.Shapes.Item(3).TextFrame.TextRange.Text = "4/15/2019"

$layouts = $sm.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shp = $layout.Shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "4/15/2019"
        }
    }
}

$nm = $p.NotesMaster
for ($j = 1; $j -le $nm.Shapes.Count; $j++) {
    $shp = $nm.Shapes.Item($j)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = "4/15/2019"
    }
}

# ---------------------------------------------------------------------
# 2) "post(TravelBuddyChangedEvent)" -> "commit()" (first post() callout,
#    :UI -> :EventsCenter lifeline region), plus reposition/resize.
# ---------------------------------------------------------------------
$shp = $s.Shapes.Item(17)   # TextBox 32 (id 33)
$shp.Left   = EMU(6894285)
$shp.Top    = EMU(2275417)
$shp.Width  = EMU(783962)
$shp.Height = EMU(215444)
$tr = $shp.TextFrame.TextRange
$tr.Characters(1, $tr.Length).Text = "commit()"

# ---------------------------------------------------------------------
# 3) ":EventsCenter" -> ":VersionedTravelBuddy" rectangle + reposition.
# ---------------------------------------------------------------------
$shp = $s.Shapes.Item(21)   # Rectangle 62 (id 39)
$shp.Left   = EMU(6858000)
$shp.Top    = EMU(1169508)
$shp.Width  = EMU(2209780)
$shp.Height = EMU(346760)
$tr = $shp.TextFrame.TextRange
$tr.Characters(2, $tr.Length - 1).Text = "VersionedTravelBuddy"

# ---------------------------------------------------------------------
# 4) Lifeline + activation bar under the renamed object, reposition only.
# ---------------------------------------------------------------------
$shp = $s.Shapes.Item(22)   # Straight Connector 39 (id 40)
$shp.Left   = EMU(7920608)
$shp.Top    = EMU(1524000)
$shp.Width  = EMU(0)
$shp.Height = EMU(1723059)

$shp = $s.Shapes.Item(23)   # Rectangle 40 (id 41)
$shp.Left   = EMU(7848600)
$shp.Top    = EMU(2540897)
$shp.Width  = EMU(142006)
$shp.Height = EMU(176787)

# ---------------------------------------------------------------------
# 5) The two arrow connectors feeding the activation bar shrink to match
#    the new (closer) lifeline position.
# ---------------------------------------------------------------------
$shp = $s.Shapes.Item(24)   # Straight Arrow Connector 41 (id 42)
$shp.Left   = EMU(5943992)
$shp.Top    = EMU(2539459)
$shp.Width  = EMU(1975611)
$shp.Height = EMU(1438)

$shp = $s.Shapes.Item(25)   # Straight Arrow Connector 43 (id 44)
$shp.Left   = EMU(5943992)
$shp.Top    = EMU(2716246)
$shp.Width  = EMU(1975611)
$shp.Height = EMU(1438)

# ---------------------------------------------------------------------
# 6) Second "post(TravelBuddyChangedEvent)" -> "commit()" callout
#    (:UI -> :EventsCenter second occurrence), plus reposition/resize.
# ---------------------------------------------------------------------
$shp = $s.Shapes.Item(29)   # TextBox 61 (id 62)
$shp.Left   = EMU(3447191)
$shp.Top    = EMU(4339076)
$shp.Width  = EMU(759234)
$shp.Height = EMU(215444)
$tr = $shp.TextFrame.TextRange
$tr.Characters(1, $tr.Length).Text = "commit()"

# ---------------------------------------------------------------------
# 7) Second ":EventsCenter" -> ":VersionedTravelBuddy" rectangle + reposition.
# ---------------------------------------------------------------------
$shp = $s.Shapes.Item(31)   # Rectangle 62 (id 66)
$shp.Left   = EMU(3429000)
$shp.Top    = EMU(3826911)
$shp.Width  = EMU(2071319)
$shp.Height = EMU(346760)
$tr = $shp.TextFrame.TextRange
$tr.Characters(2, $tr.Length - 1).Text = "VersionedTravelBuddy"

# ---------------------------------------------------------------------
# 8) "handleTravelBuddyChangedEvent()" -> "saveTravelBuddy()" + reposition.
# ---------------------------------------------------------------------
$shp = $s.Shapes.Item(37)   # TextBox 73 (id 74)
$shp.Left = EMU(5276972)
$shp.Top  = EMU(4599349)
$tr = $shp.TextFrame.TextRange
$tr.Characters(1, 29).Text = "saveTravelBuddy"

# ---------------------------------------------------------------------
# 9) "handleTravelBuddyChangedEvent()" -> "indicateModified()" + reposition.
# ---------------------------------------------------------------------
$shp = $s.Shapes.Item(44)   # TextBox 49 (id 50)
$shp.Left = EMU(1957403)
$shp.Top  = EMU(4917282)
$tr = $shp.TextFrame.TextRange
$tr.Characters(1, 29).Text = "indicateModified"
